$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (AEDB.CEA / MCP1_pg_ug_2015_rank / AsymptSympt) values ---
$ws.Range("D2").Value = 0.231510840395742
$ws.Range("E2").Value = 0.110871661398763
$ws.Range("F2").Value = 1.2605029908548
$ws.Range("G2").Value = 1.01430371706572
$ws.Range("H2").Value = 1.56646156690655
$ws.Range("I2").Value = 2.08809751270061
$ws.Range("J2").Value = 0.0367890393838166
$ws.Range("K2").Value = 0.0825603135493815
$ws.Range("L2").Value = 0.056179761222278
$ws.Range("M2").Value = 0.111560948396554
$ws.Range("N2").Value = 2423
$ws.Range("O2").Value = 1038
$ws.Range("P2").Value = 57.1605447791993

# --- New row 3: AEDB.CEA / MCP1_pg_ml_2015_rank / AsymptSympt ---
$ws.Range("A3").Value = "AEDB.CEA"
$ws.Range("B3").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C3").Value = "AsymptSympt"
$ws.Range("D3").Value = 0.363828497456171
$ws.Range("E3").Value = 0.119157758680502
$ws.Range("F3").Value = 1.43882743043218
$ws.Range("G3").Value = 1.13914649162784
$ws.Range("H3").Value = 1.81734692577223
$ws.Range("I3").Value = 3.05333451623326
$ws.Range("J3").Value = 0.002263135045014
$ws.Range("K3").Value = 0.0896426922742771
$ws.Range("L3").Value = 0.0608495285666619
$ws.Range("M3").Value = 0.120834104109508
$ws.Range("N3").Value = 2423
$ws.Range("O3").Value = 1038
$ws.Range("P3").Value = 57.1605447791993

# --- New row 4: AEDB.CEA / MCP1_rank / AsymptSympt ---
$ws.Range("A4").Value = "AEDB.CEA"
$ws.Range("B4").Value = "MCP1_rank"
$ws.Range("C4").Value = "AsymptSympt"
$ws.Range("D4").Value = 0.35843253654479
$ws.Range("E4").Value = 0.131412121172534
$ws.Range("F4").Value = 1.43108448296749
$ws.Range("G4").Value = 1.10612706009474
$ws.Range("H4").Value = 1.85150772571727
$ws.Range("I4").Value = 2.72754547561253
$ws.Range("J4").Value = 0.00638074617586061
$ws.Range("K4").Value = 0.0643929941859342
$ws.Range("L4").Value = 0.055581629821186
$ws.Range("M4").Value = 0.0944375497408653
$ws.Range("N4").Value = 2423
$ws.Range("O4").Value = 498
$ws.Range("P4").Value = 79.4469665703673
